$wb = $excel.ActiveWorkbook

# Insert a new worksheet "House part sizes" before the first existing sheet
$firstSheet = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($firstSheet)
$newSheet.Name = "House part sizes"

# Populate the new sheet with part sizes data
$newSheet.Range("A1").Value = "Part"
$newSheet.Range("B1").Value = "Square meters"
$newSheet.Range("A2").Value = "Roof"
$newSheet.Range("B2").Value = 50
$newSheet.Range("A3").Value = "External walls"
$newSheet.Range("B3").Value = 100
$newSheet.Range("A4").Value = "Foundation"
$newSheet.Range("B4").Value = 150

# Turn the data range into a table, matching the style used elsewhere
$list = $newSheet.ListObjects.Add(1, $newSheet.Range("A1:B4"), $null, 1)
$list.Name = "Table3"
$list.TableStyle = "TableStyleLight9"

$newSheet.Columns.Item(2).ColumnWidth = (13.9296875 - (5/6))

# Restore the selection on the "Supplier and cost" sheet (no longer the active tab)
$supplierSheet = $wb.Worksheets.Item("Supplier and cost")
[void]$supplierSheet.Range("F28").Select()

# Make the new sheet the active / selected tab with its own selection
$newSheet.Activate()
[void]$newSheet.Range("B20").Select()
